$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '30.398.97'
Set-TextValue 'E2' '  +1.09%  '
Set-TextValue 'D3' '1.851.35'
Set-TextValue 'E3' '  +1.10%  '
Set-TextValue 'D4' '1.002'
Set-TextValue 'E4' '  +0.15%  '
Set-TextValue 'D5' '233.50'
Set-TextValue 'E5' '  +1.99%  '
Set-TextValue 'E6' '  +0.24%  '
Set-TextValue 'D7' '0.4744'
Set-TextValue 'E7' '  +2.92%  '
Set-TextValue 'D8' '0.2750'
Set-TextValue 'E8' '  +2.59%  '
Set-TextValue 'D9' '0.06318'
Set-TextValue 'E9' '  +2.12%  '
Set-TextValue 'D10' '17.59'
Set-TextValue 'E10' '  +10.63%  '
Set-TextValue 'D11' '1.842.07'
Set-TextValue 'E11' '  +0.55%  '
Set-TextValue 'D12' '0.07459'
Set-TextValue 'E12' '  +1.57%  '
Set-TextValue 'D13' '4.958'
Set-TextValue 'E13' '  +1.61%  '
Set-TextValue 'D14' '84.56'
Set-TextValue 'E14' '  +2.57%  '
Set-TextValue 'D15' '0.6244'
Set-TextValue 'E15' '  +1.60%  '
Set-TextValue 'D16' '30.373.09'
Set-TextValue 'E16' '  +1.24%  '
Set-TextValue 'D17' '245.52'
Set-TextValue 'E17' '  +9.72%  '
Set-TextValue 'D18' '1.002'
Set-TextValue 'E18' '  +0.16%  '
Set-TextValue 'D19' '12.69'
Set-TextValue 'E19' '  +3.95%  '
Set-TextValue 'D20' '0.000007315'
Set-TextValue 'E20' '  +1.64%  '
Set-TextValue 'E21' '  +0.14%  '
Set-TextValue 'E22' '  +2.45%  '
Set-TextValue 'D23' '5.905'
Set-TextValue 'E23' '  +1.58%  '
Set-TextValue 'B24' 'Monero'
Set-TextValue 'C24' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D24' '164.78'
Set-TextValue 'E24' '  +0.09%  '
Set-TextValue 'B25' 'Cosmos'
Set-TextValue 'C25' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D25' '9.109'
Set-TextValue 'E25' '  +0.62%  '
Set-TextValue 'D26' '17.95'
Set-TextValue 'E26' '  +2.46%  '
Set-TextValue 'D27' '1.868'
Set-TextValue 'E27' '  +2.08%  '
Set-TextValue 'D28' '0.1030'
Set-TextValue 'E28' '  +1.72%  '
Set-TextValue 'D29' '1.354'
Set-TextValue 'E29' '  -0.97%  '
Set-TextValue 'D30' '4.046'
Set-TextValue 'E30' '  +0.15%  '
Set-TextValue 'D31' '3.823'
Set-TextValue 'E31' '  +2.28%  '
Set-TextValue 'D32' '0.04838'
Set-TextValue 'E32' '  +1.50%  '
Set-TextValue 'D33' '1.127'
Set-TextValue 'E33' '  +0.40%  '
Set-TextValue 'D34' '0.6963'
Set-TextValue 'E34' '  +0.56%  '
Set-TextValue 'D35' '2.696'
Set-TextValue 'E35' '  +0.46%  '
Set-TextValue 'D36' '0.01895'
Set-TextValue 'E36' '  +5.40%  '
Set-TextValue 'D37' '2.683'
Set-TextValue 'E37' '  +2.94%  '
Set-TextValue 'D38' '2.001'
Set-TextValue 'E38' '  +5.48%  '
Set-TextValue 'E39' '  -0.82%  '
Set-TextValue 'D40' '106.36'
Set-TextValue 'E40' '  +3.51%  '
Set-TextValue 'D41' '1.002'
Set-TextValue 'E41' '  +0.75%  '
Set-TextValue 'D42' '5.525'
Set-TextValue 'E42' '  +1.98%  '
Set-TextValue 'D43' '0.4047'
Set-TextValue 'E43' '  +2.22%  '
Set-TextValue 'D44' '7.166'
Set-TextValue 'E44' '  +4.69%  '
Set-TextValue 'D45' '62.90'
Set-TextValue 'E45' '  +6.94%  '
Set-TextValue 'D46' '0.1195'
Set-TextValue 'E46' '  +1.54%  '
Set-TextValue 'D47' '33.68'
Set-TextValue 'E47' '  +4.10%  '
Set-TextValue 'D48' '8.535'
Set-TextValue 'E48' '  +1.44%  '
Set-TextValue 'D49' '0.05517'
Set-TextValue 'E49' '  -0.07%  '
Set-TextValue 'D50' '1.348'
Set-TextValue 'E50' '  +0.54%  '
Set-TextValue 'D51' '0.3678'
Set-TextValue 'E51' '  +2.48%  '
